# Auto-generated: rotate rows 9-13 content per the target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel's automatic
# type inference (dates / numeric-looking strings), matching the source
# workbook where these are stored as plain text.
function Set-TextValue($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = $ws.Range("C1").Style
}

# Row 9 <= current row 13 content
$ws.Range("A9").Value = 112093595
$ws.Range("B9").Value = 90806
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 4361
$ws.Range("F9").Value = "Orange taggsvamp"
$ws.Range("G9").Value = "Hydnellum aurantiacum"
$ws.Range("H9").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I9").Value = ""
$ws.Range("K9").Value = "teleomorf"
$ws.Range("P9").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("Q9").Value = 653792
$ws.Range("R9").Value = 6576998
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Stockholm"
$ws.Range("U9").Value = "Ekerö"
$ws.Range("V9").Value = "Uppland"
$ws.Range("W9").Value = "Ekerö"
Set-TextValue $ws "Y9" "2023-09-14"
Set-TextValue $ws "AA9" "2023-09-14"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AI9").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = "Jan Yngve Andersson"
$ws.Range("AX9").Value = "Jan Yngve Andersson"
$ws.Range("AY9").Value = ""

# Row 10 <= current row 12 content
$ws.Range("A10").Value = 112093593
$ws.Range("B10").Value = 103781
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 221144
$ws.Range("F10").Value = "Grönpyrola"
$ws.Range("G10").Value = "Pyrola chlorantha"
$ws.Range("H10").Value = "Sw."
$ws.Range("I10").Value = ""
$ws.Range("K10").Value = "fullt utvecklade blad"
$ws.Range("P10").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("Q10").Value = 653786
$ws.Range("R10").Value = 6577035
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Stockholm"
$ws.Range("U10").Value = "Ekerö"
$ws.Range("V10").Value = "Uppland"
$ws.Range("W10").Value = "Ekerö"
Set-TextValue $ws "Y10" "2023-09-14"
Set-TextValue $ws "AA10" "2023-09-14"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AI10").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AT10").Value = ""
$ws.Range("AW10").Value = "Jan Yngve Andersson"
$ws.Range("AX10").Value = "Jan Yngve Andersson"
$ws.Range("AY10").Value = ""

# Row 11 <= current row 9 content
$ws.Range("A11").Value = 112097135
$ws.Range("B11").Value = 90814
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 4364
$ws.Range("F11").Value = "Dropptaggsvamp"
$ws.Range("G11").Value = "Hydnellum ferrugineum"
$ws.Range("H11").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = "teleomorf"
$ws.Range("N11").Value = ""
$ws.Range("P11").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("Q11").Value = 653889
$ws.Range("R11").Value = 6576888
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = "Stockholm"
$ws.Range("U11").Value = "Ekerö"
$ws.Range("V11").Value = "Uppland"
$ws.Range("W11").Value = "Ekerö"
Set-TextValue $ws "Y11" "2023-09-14"
Set-TextValue $ws "AA11" "2023-09-14"
$ws.Range("AC11").Value = "Mörkröda droppar på hattöversidan. Smak besk efter ett långt tag, ej brännande. Köttet färgas mörkviolett med KOH."
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AF11").Value = ""
$ws.Range("AG11").Value = $false
$ws.Range("AI11").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AT11").Value = ""
$ws.Range("AW11").Value = "Jan Yngve Andersson"
$ws.Range("AX11").Value = "Jan Yngve Andersson"
$ws.Range("AY11").Value = ""

# Row 12 <= current row 11 content
$ws.Range("A12").Value = 112093592
$ws.Range("B12").Value = 90806
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = "Orange taggsvamp"
$ws.Range("G12").Value = "Hydnellum aurantiacum"
$ws.Range("H12").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I12").Value = ""
$ws.Range("K12").Value = "teleomorf"
$ws.Range("P12").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("Q12").Value = 653789
$ws.Range("R12").Value = 6577029
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = "Stockholm"
$ws.Range("U12").Value = "Ekerö"
$ws.Range("V12").Value = "Uppland"
$ws.Range("W12").Value = "Ekerö"
Set-TextValue $ws "Y12" "2023-09-14"
Set-TextValue $ws "AA12" "2023-09-14"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AI12").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AT12").Value = ""
$ws.Range("AW12").Value = "Jan Yngve Andersson"
$ws.Range("AX12").Value = "Jan Yngve Andersson"
$ws.Range("AY12").Value = ""

# Row 13 <= current row 10 content
$ws.Range("A13").Value = 112093591
$ws.Range("B13").Value = 103742
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "EN"
$ws.Range("E13").Value = 340
$ws.Range("F13").Value = "Ryl"
$ws.Range("G13").Value = "Chimaphila umbellata"
$ws.Range("H13").Value = "(L.) W. P. C. Barton"
Set-TextValue $ws "I13" "9"
$ws.Range("J13").Value = "stjälkar/strån/skott"
$ws.Range("K13").Value = "fullt utvecklade blad"
$ws.Range("P13").Value = "Svarvartorp ca 400 m SO om, Upl"
$ws.Range("Q13").Value = 653798
$ws.Range("R13").Value = 6576988
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = "Stockholm"
$ws.Range("U13").Value = "Ekerö"
$ws.Range("V13").Value = "Uppland"
$ws.Range("W13").Value = "Ekerö"
Set-TextValue $ws "Y13" "2023-09-14"
Set-TextValue $ws "AA13" "2023-09-14"
$ws.Range("AC13").Value = "Ganska små och taniga. Översiktlig räkning och troligen finns fler på de båda dellokalerna."
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AI13").Value = "Gles barrskog på sand (både tall och gran)"
$ws.Range("AT13").Value = ""
$ws.Range("AW13").Value = "Jan Yngve Andersson"
$ws.Range("AX13").Value = "Jan Yngve Andersson"
$ws.Range("AY13").Value = ""
